$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting existing rows 210:256 down to 211:257
$ws.Rows("210:210").Insert()

# Populate the newly inserted row 210 with the new weekly record
$ws.Range("A210").Value = 10
$ws.Range("B210").Value = "Vega Modelo de Temuco"
$ws.Range("C210").Value = "La Araucanía"
$ws.Range("D210").Value = 44932
$ws.Range("E210").Value = 9
$ws.Range("F210").Value = 100112005
$ws.Range("G210").Value = "Puerro"
$ws.Range("H210").Value = "Azul de Maquehue"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 30
$ws.Range("K210").Value = 15000
$ws.Range("L210").Value = 15000
$ws.Range("M210").Value = 15000
$ws.Range("N210").Value = "$/docena de paquetes"
$ws.Range("O210").Value = "Provincia de Cautín"
$ws.Range("P210").Value = 1250
$ws.Range("Q210").Value = 12
$ws.Range("R210").Value = "Hortaliza"

# Apply the same date number format used by the other rows in column D
$ws.Range("D210").NumberFormat = $ws.Range("D211").NumberFormat
